$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of cell address -> new display value (kept as text to preserve
# the original "dotted" price formatting, e.g. "27.851.33").
$updates = @{
    "D2" = "27.851.33"
    "E2" = "  +0.57%  "
    "D3" = "1.751.85"
    "E3" = "  -0.69%  "
    "E4" = "  +0.26%  "
    "D5" = "333.70"
    "E5" = "  -0.27%  "
    "D6" = "1.002"
    "E6" = "  +0.13%  "
    "D7" = "0.3871"
    "E7" = "  +1.31%  "
    "D8" = "0.3385"
    "E8" = "  -1.07%  "
    "D9" = "45.61"
    "E9" = "  -2.28%  "
    "D10" = "1.113"
    "E10" = "  -1.51%  "
    "E11" = "  -2.41%  "
    "E12" = "  +0.34%  "
    "D13" = "22.49"
    "E13" = "  +0.96%  "
    "D14" = "6.180"
    "E14" = "  -2.34%  "
    "D15" = "1.754.11"
    "E15" = "  -0.40%  "
    "D16" = "7.094"
    "E16" = "  +0.60%  "
    "E17" = "  -1.13%  "
    "E18" = "  -1.05%  "
    "D19" = "79.42"
    "E19" = "  -2.75%  "
    "D20" = "1.001"
    "E20" = "  +0.00%  "
    "E21" = "  -2.66%  "
    "D22" = "6.180"
    "E22" = "  -3.13%  "
    "D23" = "27.859.60"
    "E23" = "  +0.70%  "
    "E24" = "  -2.73%  "
    "D25" = "2.400"
    "E25" = "  +0.79%  "
    "D26" = "154.69"
    "E26" = "  +1.39%  "
    "D27" = "19.85"
    "E27" = "  -3.70%  "
    "D28" = "2.300"
    "E28" = "  -4.16%  "
    "D29" = "1.951.13"
    "E29" = "  -0.58%  "
    "D30" = "1.293"
    "E30" = "  -9.22%  "
    "D31" = "130.49"
    "E31" = "  -2.69%  "
    "D32" = "4.019"
    "E32" = "  +1.55%  "
    "D33" = "5.811"
    "E33" = "  -4.30%  "
    "D34" = "0.08789"
    "E34" = "  +0.57%  "
    "D35" = "12.16"
    "E35" = "  -4.04%  "
    "D36" = "1.540"
    "E36" = "  +2.86%  "
    "D37" = "0.6536"
    "E37" = "  -3.30%  "
    "D38" = "5.133"
    "E38" = "  -2.81%  "
    "D39" = "0.02276"
    "E39" = "  -5.54%  "
    "D40" = "0.06111"
    "E40" = "  -2.52%  "
    "D41" = "0.2107"
    "E41" = "  -2.93%  "
    "D42" = "1.207"
    "E42" = "  -3.45%  "
    "D43" = "8.024"
    "E43" = "  -2.04%  "
    "E44" = "  +0.05%  "
    "D45" = "13.75"
    "E45" = "  -2.26%  "
    "D46" = "3.819"
    "E46" = "  -0.50%  "
    "D47" = "0.6037"
    "E47" = "  -3.16%  "
    "D48" = "126.91"
    "E48" = "  -3.22%  "
    "D49" = "1.993"
    "E49" = "  -3.23%  "
    "D50" = "1.115"
    "E50" = "  +5.20%  "
    "D51" = "1.160"
    "E51" = "  +1.72%  "
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Force text storage so numeric-looking strings (e.g. "333.70",
    # "1.540") are not reinterpreted as numbers and lose their
    # trailing/format-significant digits.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    # Restore the default (unstyled) look - matches the workbook's
    # original unstyled data cells.
    $cell.Style = "Normal"
}
